$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Create two brand-new rows (114, 115) by copying the format of row 113 ---
# (row 113 currently ends the sheet and has blank C/D cells, which is exactly
# the pattern we want for the two new trailing rows as well).
$ws.Range("A113:D113").Copy()
$ws.Range("A114:D115").PasteSpecial(-4122)

# --- Column C/D moving-sum values for rows 90-91 (recomputed after the new
# data point inserted below changes the 7-day window) ---
$ws.Range("C90").Value = 54
$ws.Range("D90").Value = 1022.727272727273
$ws.Range("C91").Value = 53
$ws.Range("D91").Value = 1003.787878787879

# row 92 (44234) is unchanged

# --- A new observation (date 44235) was inserted, shifting every
# subsequent row down by one and changing the recomputed moving sums ---
$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = 50
$ws.Range("D93").Value = 946.969696969697

$ws.Range("A94").Value = 44236
$ws.Range("B94").Value = 0
$ws.Range("C94").Value = 44
$ws.Range("D94").Value = 833.3333333333334

$ws.Range("A95").Value = 44237
$ws.Range("B95").Value = 1
$ws.Range("C95").Value = 36
$ws.Range("D95").Value = 681.8181818181818

$ws.Range("A96").Value = 44238
$ws.Range("B96").Value = 2
$ws.Range("C96").Value = 22
$ws.Range("D96").Value = 416.6666666666667

$ws.Range("A97").Value = 44239
$ws.Range("B97").Value = 4
$ws.Range("C97").Value = 30
$ws.Range("D97").Value = 568.1818181818182

$ws.Range("A98").Value = 44240
$ws.Range("B98").Value = 2
$ws.Range("C98").Value = 30
$ws.Range("D98").Value = 568.1818181818182

$ws.Range("A99").Value = 44241
$ws.Range("B99").Value = 11
$ws.Range("C99").Value = 30
$ws.Range("D99").Value = 568.1818181818182

$ws.Range("A100").Value = 44242
$ws.Range("B100").Value = 10
$ws.Range("C100").Value = 29
$ws.Range("D100").Value = 549.2424242424242

$ws.Range("A101").Value = 44243
$ws.Range("B101").Value = 0
$ws.Range("C101").Value = 32
$ws.Range("D101").Value = 606.060606060606

$ws.Range("A102").Value = 44244
$ws.Range("B102").Value = 1
$ws.Range("C102").Value = 34
$ws.Range("D102").Value = 643.9393939393939

$ws.Range("A103").Value = 44245
$ws.Range("B103").Value = 1
$ws.Range("C103").Value = 26
$ws.Range("D103").Value = 492.4242424242424

$ws.Range("A104").Value = 44246
$ws.Range("B104").Value = 7
$ws.Range("C104").Value = 21
$ws.Range("D104").Value = 397.7272727272727

$ws.Range("A105").Value = 44247
$ws.Range("B105").Value = 4
$ws.Range("C105").Value = 26
$ws.Range("D105").Value = 492.4242424242424

$ws.Range("A106").Value = 44248
$ws.Range("B106").Value = 3
$ws.Range("C106").Value = 26
$ws.Range("D106").Value = 492.4242424242424

$ws.Range("A107").Value = 44249
$ws.Range("B107").Value = 5
$ws.Range("C107").Value = 31
$ws.Range("D107").Value = 587.1212121212121

$ws.Range("A108").Value = 44250
$ws.Range("B108").Value = 5
$ws.Range("C108").Value = 32
$ws.Range("D108").Value = 606.060606060606

$ws.Range("A109").Value = 44251
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = 33
$ws.Range("D109").Value = 625

$ws.Range("A110").Value = 44252
$ws.Range("B110").Value = 6
$ws.Range("C110").Value = 34
$ws.Range("D110").Value = 643.9393939393939

$ws.Range("A111").Value = 44253
$ws.Range("B111").Value = 8
$ws.Range("C111").Value = 45
$ws.Range("D111").Value = 852.2727272727273

$ws.Range("A112").Value = 44254
$ws.Range("B112").Value = 5
$ws.Range("C112").Value = 41
$ws.Range("D112").Value = 776.5151515151515

# row 113 keeps its blank C/D (copied format already matches); only A/B change
$ws.Range("A113").Value = 44255
$ws.Range("B113").Value = 4

# rows 114-115 are the two brand-new trailing rows (format already copied
# above from old row 113, so C/D remain blank)
$ws.Range("A114").Value = 44256
$ws.Range("B114").Value = 16

$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 1
